$d = $word.ActiveDocument

# Fix the "-s pandoc.css" typo to "-c pandoc.css" in the three repeated
# pandoc command-line examples.
$d.Content.Find.Execute("pandoc -s -S -s pandoc.css --bibliography=library.bib", $true, $false, $false, $false, $false, $true, 1, $false, "pandoc -s -S -c pandoc.css --bibliography=library.bib", 2)

# Drop the "docs/" prefix from the output file paths in the same examples.
$d.Content.Find.Execute("--csl=mee.csl test.md -o docs/test.html", $true, $false, $false, $false, $false, $true, 1, $false, "--csl=mee.csl test.md -o test.html", 2)
$d.Content.Find.Execute("--csl=mee.csl test.md -o docs/test.docx", $true, $false, $false, $false, $false, $true, 1, $false, "--csl=mee.csl test.md -o test.docx", 2)
$d.Content.Find.Execute("--csl=mee.csl test.md -o docs/test.pdf", $true, $false, $false, $false, $false, $true, 1, $false, "--csl=mee.csl test.md -o test.pdf", 2)
